# "pre excel sheet update"
# The workbook originally has three sheets: "sector3", "sector1" and "Sheet".
# The first two each host a picture (via a drawing) anchored in cell A1.
# Only the (already empty) "Sheet" worksheet should remain afterwards.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the sheets (and their embedded picture / drawing) that are not
# being kept, leaving "Sheet" as the sole, active worksheet.
$wb.Worksheets.Item("sector3").Delete()
$wb.Worksheets.Item("sector1").Delete()

# Make sure the surviving sheet is selected/active.
$wb.Worksheets.Item("Sheet").Activate()
